$d = $word.ActiveDocument

# Locate the target paragraph: "Yes, t" + "he solution to draw a hand ..."
# -- this is the one that also carries the stray _GoBack bookmark sitting
# mid-sentence, between the two runs.
$target = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "*rhythm of the counting.*") {
        $target = $p
        break
    }
}
if ($target -eq $null) {
    throw "Could not locate the target paragraph (rhythm of the counting)."
}

$w = "xmlns:w='http://schemas.openxmlformats.org/wordprocessingml/2006/main'"
$w14 = "xmlns:w14='http://schemas.microsoft.com/office/word/2010/wordml'"
$ns = "$w $w14"

# Paragraph 1 keeps its original identity/attrs, the bookmark is gone from
# the middle of the sentence, and the two runs remain split as they were.
$para1 = "<w:p $ns w14:paraId='468C5088' w14:textId='426AE5EB' w:rsidR='00A42C06' w:rsidRDefault='00F02B18' w:rsidP='00A42C06'>" +
         "<w:pPr><w:pStyle w:val='ListParagraph'/><w:numPr><w:ilvl w:val='1'/><w:numId w:val='7'/></w:numPr></w:pPr>" +
         "<w:r><w:t>Yes, t</w:t></w:r>" +
         "<w:r w:rsidR='00A42C06'><w:t>he solution to draw a hand and label it with numbers and figuring out a system helps in figuring out the rhythm of the counting.</w:t></w:r>" +
         "</w:p>"

# Brand new list item (same numbered-list level) that continues the answer
# and now carries the relocated _GoBack bookmark at its very end.
$para2 = "<w:p $ns>" +
         "<w:pPr><w:pStyle w:val='ListParagraph'/><w:numPr><w:ilvl w:val='1'/><w:numId w:val='7'/></w:numPr></w:pPr>" +
         "<w:r><w:t>Yes, because the visual aid will help to figure out the system quickly.</w:t></w:r>" +
         "<w:bookmarkStart w:id='0' w:name='_GoBack'/><w:bookmarkEnd w:id='0'/>" +
         "</w:p>"

[void]$target.Range.InsertXML($para1 + $para2)
